# Apply updated "dSF" (column F) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -2
    17 = 8
    22 = 4
    26 = -1
    27 = 4
    28 = 6
    31 = -5
    33 = -5
    34 = 2
    36 = -6
    37 = -5
    40 = -4
    46 = -3
    47 = 6
    48 = 7
    50 = -3
    51 = -1
    52 = -6
    53 = -5
    54 = -8
    55 = -8
    57 = 8
    58 = -5
    60 = -3
    61 = -1
    62 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
